# The workbook has two sheets: "Affichage" and "Candidatures".
# The edit adds one new row (row 5) to the "Candidatures" sheet, which is a
# duplicate of row 4 ("Charlie C" / PHY2710, PHY2400 / 3.42 / Plasmas) except
# the name in column A becomes "Denise D" (a new entry, i.e. "adding 1
# equality" to the example).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Candidatures")
$ws.Activate() | Out-Null

# Duplicate row 4 (values + formatting) into the new row 5 ...
$ws.Range("A4:H4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial() | Out-Null

# ... then overwrite the name column with the new "Denise D" entry.
$ws.Range("A5").Value2 = "Denise D"

# Reflect the post-edit cursor position/selection (next empty row).
$ws.Range("A6").Select() | Out-Null
